# issue #5: stock data from json to db
#
# The stock ("股票") sheet gains three new columns that were introduced when
# the source data moved from ad-hoc JSON scraping to a proper database
# import: `category` (inserted right after property_category/before date),
# and `source_file` + `index` (appended after legislator_id), recording
# which import batch / row each record came from.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert the new "category" column before the existing "date" column,
# shifting date / legislator_name / legislator_id one column to the right
# (I -> J -> K -> L).
$ws.Columns.Item(9).Insert()

# Append two more new columns after legislator_id (L) for "source_file"
# and "index".
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(14).Insert()

# "category" column (I)
$ws.Range("I1").Value = "category"
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("I4").Value = "normal"
$ws.Range("I5").Value = "normal"

# "source_file" column (M)
$ws.Range("M1").Value = "source_file"
$ws.Range("M2").Value = "tmp2691"
$ws.Range("M3").Value = "tmp2691"
$ws.Range("M4").Value = "tmp2691"
$ws.Range("M5").Value = "tmp2691"

# "index" column (N) - the original row index from the source data
$ws.Range("N1").Value = "index"
$ws.Range("N2").Value = 58
$ws.Range("N3").Value = 59
$ws.Range("N4").Value = 60
$ws.Range("N5").Value = 61
